$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.863.15'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '3.329.30'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'582.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = "'177.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").Value = '3.324.82'
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("E10").Value = '  +4.67%  '
$ws.Range("D11").Value = "'0.583"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = "'47.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.54%  '
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").Value = "'700.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("D15").Value = '3.868.50'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = "'8.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '67.906.80'
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '3.336.93'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = "'17.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = "'11.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").Value = "'5.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.83%  '
$ws.Range("D24").Value = "'17.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = "'99.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = "'3.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = "'9.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("D29").Value = "'33.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("D30").Value = "'8.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("D31").Value = "'7.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.00%  '
$ws.Range("D32").Value = "'568.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("D33").Value = "'11.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("E34").Value = '  +2.44%  '
$ws.Range("D35").Value = "'57.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.30%  '
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").Value = '3.687.40'
$ws.Range("E37").Value = '  -5.25%  '
$ws.Range("D38").Value = "'3.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("D39").Value = "'34.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.21%  '
$ws.Range("E40").Value = '  +2.92%  '
$ws.Range("E41").Value = '  +2.21%  '
$ws.Range("E42").Value = '  +5.79%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0676'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = "'0.337"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = "'3.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.68%  '
$ws.Range("D46").Value = "'0.0407"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = "'2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.37%  '
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = "'1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.25%  '
$ws.Range("D51").Value = "'130.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.40%  '
